$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 128.4548946666667
$ws.Range("H2").Value = 385.364684
$ws.Range("I2").Value = 0.2815548034715028
$ws.Range("J2").Value = 0.2815548034715028
$ws.Range("M2").Value = 1.646588666666666
$ws.Range("N2").Value = 4.939766
$ws.Range("O2").Value = 0.039310317935267
$ws.Range("P2").Value = 0.039310317935267
$ws.Range("Q2").Value = 211.5123737359937
$ws.Range("R2").Value = 1903.611363623944
$ws.Range("S2").Value = 0.01106800884066639
$ws.Range("T2").Value = 0.01106800884066639
$ws.Range("G3").Value = 128.4548946666667
$ws.Range("H3").Value = 385.364684
$ws.Range("I3").Value = 0.2815548034715028
$ws.Range("J3").Value = 0.2815548034715028
$ws.Range("O3").Value = 0.278787195370394
$ws.Range("P3").Value = 0.278787195370394
$ws.Range("Q3").Value = 1500.037256302383
$ws.Range("R3").Value = 13500.33530672145
$ws.Range("S3").Value = 0.07849387400288273
$ws.Range("T3").Value = 0.07849387400288273
$ws.Range("G4").Value = 128.4548946666667
$ws.Range("H4").Value = 385.364684
$ws.Range("I4").Value = 0.2815548034715028
$ws.Range("J4").Value = 0.2815548034715028
$ws.Range("M4").Value = 0.7553226666666667
$ws.Range("N4").Value = 2.265968
$ws.Range("O4").Value = 0.01803241742850595
$ws.Range("P4").Value = 0.01803241742850595
$ws.Range("Q4").Value = 97.02489358601244
$ws.Range("R4").Value = 873.224042274112
$ws.Range("S4").Value = 0.005077113745199094
$ws.Range("T4").Value = 0.005077113745199094
$ws.Range("G5").Value = 128.4548946666667
$ws.Range("H5").Value = 385.364684
$ws.Range("I5").Value = 0.2815548034715028
$ws.Range("J5").Value = 0.2815548034715028
$ws.Range("M5").Value = 27.21325766666666
$ws.Range("N5").Value = 81.63977299999999
$ws.Range("O5").Value = 0.6496836961088899
$ws.Range("P5").Value = 0.6496836961088899
$ws.Range("Q5").Value = 3495.676147108525
$ws.Range("R5").Value = 31461.08532397673
$ws.Range("S5").Value = 0.182921565376578
$ws.Range("T5").Value = 0.182921565376578
$ws.Range("G6").Value = 128.4548946666667
$ws.Range("H6").Value = 385.364684
$ws.Range("I6").Value = 0.2815548034715028
$ws.Range("J6").Value = 0.2815548034715028
$ws.Range("M6").Value = 0.5942236666666667
$ws.Range("N6").Value = 1.782671
$ws.Range("O6").Value = 0.01418637315694314
$ws.Range("P6").Value = 0.01418637315694314
$ws.Range("Q6").Value = 76.33093851010712
$ws.Range("R6").Value = 686.978446590964
$ws.Range("S6").Value = 0.003994241506176528
$ws.Range("T6").Value = 0.003994241506176528
$ws.Range("H7").Value = 457.183265
$ws.Range("I7").Value = 0.3340268313936494
$ws.Range("J7").Value = 0.3340268313936494
$ws.Range("M7").Value = 1.646588666666666
$ws.Range("N7").Value = 4.939766
$ws.Range("O7").Value = 0.039310317935267
$ws.Range("P7").Value = 0.039310317935267
$ws.Range("Q7").Value = 250.9309275795544
$ws.Range("R7").Value = 2258.37834821599
$ws.Range("S7").Value = 0.01313070094099418
$ws.Range("T7").Value = 0.01313070094099418
$ws.Range("H8").Value = 457.183265
$ws.Range("I8").Value = 0.3340268313936494
$ws.Range("J8").Value = 0.3340268313936494
$ws.Range("O8").Value = 0.278787195370394
$ws.Range("P8").Value = 0.278787195370394
$ws.Range("S8").Value = 0.09312240350269499
$ws.Range("T8").Value = 0.09312240350269499
$ws.Range("H9").Value = 457.183265
$ws.Range("I9").Value = 0.3340268313936494
$ws.Range("J9").Value = 0.3340268313936494
$ws.Range("M9").Value = 0.7553226666666667
$ws.Range("N9").Value = 2.265968
$ws.Range("O9").Value = 0.01803241742850595
$ws.Range("P9").Value = 0.01803241742850595
$ws.Range("Q9").Value = 115.1069609583911
$ws.Range("R9").Value = 1035.96264862552
$ws.Range("S9").Value = 0.006023311256011461
$ws.Range("T9").Value = 0.006023311256011461
$ws.Range("H10").Value = 457.183265
$ws.Range("I10").Value = 0.3340268313936494
$ws.Range("J10").Value = 0.3340268313936494
$ws.Range("M10").Value = 27.21325766666666
$ws.Range("N10").Value = 81.63977299999999
$ws.Range("O10").Value = 0.6496836961088899
$ws.Range("P10").Value = 0.6496836961088899
$ws.Range("Q10").Value = 4147.148663777649
$ws.Range("R10").Value = 37324.33797399884
$ws.Range("S10").Value = 0.2170117864193671
$ws.Range("T10").Value = 0.2170117864193671
$ws.Range("H11").Value = 457.183265
$ws.Range("I11").Value = 0.3340268313936494
$ws.Range("J11").Value = 0.3340268313936494
$ws.Range("M11").Value = 0.5942236666666667
$ws.Range("N11").Value = 1.782671
$ws.Range("O11").Value = 0.01418637315694314
$ws.Range("P11").Value = 0.01418637315694314
$ws.Range("Q11").Value = 90.55637202231279
$ws.Range("R11").Value = 815.007348200815
$ws.Range("S11").Value = 0.00473862927458164
$ws.Range("T11").Value = 0.00473862927458164
$ws.Range("G12").Value = 70.798157
$ws.Range("H12").Value = 212.394471
$ws.Range("I12").Value = 0.1551794599342134
$ws.Range("J12").Value = 0.1551794599342134
$ws.Range("M12").Value = 1.646588666666666
$ws.Range("N12").Value = 4.939766
$ws.Range("O12").Value = 0.039310317935267
$ws.Range("P12").Value = 0.039310317935267
$ws.Range("Q12").Value = 116.5754429370873
$ws.Range("R12").Value = 1049.178986433786
$ws.Range("S12").Value = 0.006100153907036956
$ws.Range("T12").Value = 0.006100153907036956
$ws.Range("G13").Value = 70.798157
$ws.Range("H13").Value = 212.394471
$ws.Range("I13").Value = 0.1551794599342134
$ws.Range("J13").Value = 0.1551794599342134
$ws.Range("O13").Value = 0.278787195370394
$ws.Range("P13").Value = 0.278787195370394
$ws.Range("Q13").Value = 826.7483574925515
$ws.Range("R13").Value = 7440.735217432964
$ws.Range("S13").Value = 0.04326204641415177
$ws.Range("T13").Value = 0.04326204641415177
$ws.Range("G14").Value = 70.798157
$ws.Range("H14").Value = 212.394471
$ws.Range("I14").Value = 0.1551794599342134
$ws.Range("J14").Value = 0.1551794599342134
$ws.Range("M14").Value = 0.7553226666666667
$ws.Range("N14").Value = 2.265968
$ws.Range("O14").Value = 0.01803241742850595
$ws.Range("P14").Value = 0.01803241742850595
$ws.Range("Q14").Value = 53.47545274032534
$ws.Range("R14").Value = 481.279074662928
$ws.Range("S14").Value = 0.00279826079786385
$ws.Range("T14").Value = 0.00279826079786385
$ws.Range("G15").Value = 70.798157
$ws.Range("H15").Value = 212.394471
$ws.Range("I15").Value = 0.1551794599342134
$ws.Range("J15").Value = 0.1551794599342134
$ws.Range("M15").Value = 27.21325766666666
$ws.Range("N15").Value = 81.63977299999999
$ws.Range("O15").Value = 0.6496836961088899
$ws.Range("P15").Value = 0.6496836961088899
$ws.Range("Q15").Value = 1926.64848876612
$ws.Range("R15").Value = 17339.83639889508
$ws.Range("S15").Value = 0.1008175650902411
$ws.Range("T15").Value = 0.1008175650902411
$ws.Range("G16").Value = 70.798157
$ws.Range("H16").Value = 212.394471
$ws.Range("I16").Value = 0.1551794599342134
$ws.Range("J16").Value = 0.1551794599342134
$ws.Range("M16").Value = 0.5942236666666667
$ws.Range("N16").Value = 1.782671
$ws.Range("O16").Value = 0.01418637315694314
$ws.Range("P16").Value = 0.01418637315694314
$ws.Range("Q16").Value = 42.06994044578234
$ws.Range("R16").Value = 378.6294640120411
$ws.Range("S16").Value = 0.002201433724919658
$ws.Range("T16").Value = 0.002201433724919658
$ws.Range("G17").Value = 20.703408
$ws.Range("H17").Value = 62.110224
$ws.Range("I17").Value = 0.04537891674549766
$ws.Range("J17").Value = 0.04537891674549767
$ws.Range("M17").Value = 1.646588666666666
$ws.Range("N17").Value = 4.939766
$ws.Range("O17").Value = 0.039310317935267
$ws.Range("P17").Value = 0.039310317935267
$ws.Range("Q17").Value = 34.089996974176
$ws.Range("R17").Value = 306.809972767584
$ws.Range("S17").Value = 0.001783859644823525
$ws.Range("T17").Value = 0.001783859644823525
$ws.Range("G18").Value = 20.703408
$ws.Range("H18").Value = 62.110224
$ws.Range("I18").Value = 0.04537891674549766
$ws.Range("J18").Value = 0.04537891674549767
$ws.Range("O18").Value = 0.278787195370394
$ws.Range("P18").Value = 0.278787195370394
$ws.Range("Q18").Value = 241.764888858592
$ws.Range("R18").Value = 2175.883999727329
$ws.Range("S18").Value = 0.0126510609284239
$ws.Range("T18").Value = 0.0126510609284239
$ws.Range("G19").Value = 20.703408
$ws.Range("H19").Value = 62.110224
$ws.Range("I19").Value = 0.04537891674549766
$ws.Range("J19").Value = 0.04537891674549767
$ws.Range("M19").Value = 0.7553226666666667
$ws.Range("N19").Value = 2.265968
$ws.Range("O19").Value = 0.01803241742850595
$ws.Range("P19").Value = 0.01803241742850595
$ws.Range("Q19").Value = 15.637753339648
$ws.Range("R19").Value = 140.739780056832
$ws.Range("S19").Value = 0.0008182915692082325
$ws.Range("T19").Value = 0.0008182915692082326
$ws.Range("G20").Value = 20.703408
$ws.Range("H20").Value = 62.110224
$ws.Range("I20").Value = 0.04537891674549766
$ws.Range("J20").Value = 0.04537891674549767
$ws.Range("M20").Value = 27.21325766666666
$ws.Range("N20").Value = 81.63977299999999
$ws.Range("O20").Value = 0.6496836961088899
$ws.Range("P20").Value = 0.6496836961088899
$ws.Range("Q20").Value = 563.407176482128
$ws.Range("R20").Value = 5070.664588339152
$ws.Range("S20").Value = 0.02948194235663252
$ws.Range("T20").Value = 0.02948194235663252
$ws.Range("G21").Value = 20.703408
$ws.Range("H21").Value = 62.110224
$ws.Range("I21").Value = 0.04537891674549766
$ws.Range("J21").Value = 0.04537891674549767
$ws.Range("M21").Value = 0.5942236666666667
$ws.Range("N21").Value = 1.782671
$ws.Range("O21").Value = 0.01418637315694314
$ws.Range("P21").Value = 0.01418637315694314
$ws.Range("Q21").Value = 12.302455014256
$ws.Range("R21").Value = 110.722095128304
$ws.Range("S21").Value = 0.0006437622464094855
$ws.Range("T21").Value = 0.0006437622464094856
$ws.Range("G22").Value = 83.88319133333333
$ws.Range("H22").Value = 251.649574
$ws.Range("I22").Value = 0.1838599884551367
$ws.Range("J22").Value = 0.1838599884551367
$ws.Range("M22").Value = 1.646588666666666
$ws.Range("N22").Value = 4.939766
$ws.Range("O22").Value = 0.039310317935267
$ws.Range("P22").Value = 0.039310317935267
$ws.Range("Q22").Value = 138.1211121732982
$ws.Range("R22").Value = 1243.090009559684
$ws.Range("S22").Value = 0.007227594601745943
$ws.Range("T22").Value = 0.007227594601745943
$ws.Range("G23").Value = 83.88319133333333
$ws.Range("H23").Value = 251.649574
$ws.Range("I23").Value = 0.1838599884551367
$ws.Range("J23").Value = 0.1838599884551367
$ws.Range("O23").Value = 0.278787195370394
$ws.Range("P23").Value = 0.278787195370394
$ws.Range("Q23").Value = 979.5493780447808
$ws.Range("R23").Value = 8815.944402403029
$ws.Range("S23").Value = 0.05125781052224057
$ws.Range("T23").Value = 0.05125781052224057
$ws.Range("G24").Value = 83.88319133333333
$ws.Range("H24").Value = 251.649574
$ws.Range("I24").Value = 0.1838599884551367
$ws.Range("J24").Value = 0.1838599884551367
$ws.Range("M24").Value = 0.7553226666666667
$ws.Range("N24").Value = 2.265968
$ws.Range("O24").Value = 0.01803241742850595
$ws.Range("P24").Value = 0.01803241742850595
$ws.Range("Q24").Value = 63.35887576640356
$ws.Range("R24").Value = 570.229881897632
$ws.Range("S24").Value = 0.003315440060223309
$ws.Range("T24").Value = 0.003315440060223309
$ws.Range("G25").Value = 83.88319133333333
$ws.Range("H25").Value = 251.649574
$ws.Range("I25").Value = 0.1838599884551367
$ws.Range("J25").Value = 0.1838599884551367
$ws.Range("M25").Value = 27.21325766666666
$ws.Range("N25").Value = 81.63977299999999
$ws.Range("O25").Value = 0.6496836961088899
$ws.Range("P25").Value = 0.6496836961088899
$ws.Range("Q25").Value = 2282.7348996563
$ws.Range("R25").Value = 20544.6140969067
$ws.Range("S25").Value = 0.119450836866071
$ws.Range("T25").Value = 0.119450836866071
$ws.Range("G26").Value = 83.88319133333333
$ws.Range("H26").Value = 251.649574
$ws.Range("I26").Value = 0.1838599884551367
$ws.Range("J26").Value = 0.1838599884551367
$ws.Range("M26").Value = 0.5942236666666667
$ws.Range("N26").Value = 1.782671
$ws.Range("O26").Value = 0.01418637315694314
$ws.Range("P26").Value = 0.01418637315694314
$ws.Range("Q26").Value = 49.84537752579489
$ws.Range("R26").Value = 448.608397732154
$ws.Range("S26").Value = 0.002608306404855826
$ws.Range("T26").Value = 0.002608306404855826
